# "build a apresentar ao negocio"
#
# This script reproduces, via the Excel COM object model, the changes that
# were captured in the OOXML diff:
#   1. Rename sheet "label_map" -> "LabelMap"
#   2. Replace the 5 shared "\\Folder\Ficheiro N.pdf" placeholder strings
#      (referenced from IDTemplates!C3,C4,C7,C8,C9,C11,C13,C14) with the
#      real local docx template paths.
#   3. Swap columns D and E (and their header/content) on IdentifEntidade.
#   4. Update the selected cell / active sheet view state on several sheets,
#      moving the "active" tab from IdentifEntidade to LabelMap.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the "label_map" worksheet to "LabelMap"
# ---------------------------------------------------------------------
$wsLabelMap = $wb.Worksheets.Item("label_map")
$wsLabelMap.Name = "LabelMap"

# ---------------------------------------------------------------------
# 2. IDTemplates: point the template-file column at the new local paths
# ---------------------------------------------------------------------
$wsTemplates = $wb.Worksheets.Item("IDTemplates")

$wsTemplates.Range("C3").Value  = 'C:\Users\brunofilipe.lobo\OneDrive - CGI\Code\realvidaseguros\Anexos\Ficheiro1.docx;'
$wsTemplates.Range("C4").Value  = 'C:\Users\brunofilipe.lobo\OneDrive - CGI\Code\realvidaseguros\Anexos\Ficheiro1.docx;'
$wsTemplates.Range("C7").Value  = 'C:\Users\brunofilipe.lobo\OneDrive - CGI\Code\realvidaseguros\Anexos\Ficheiro1.docx;'
$wsTemplates.Range("C8").Value  = 'C:\Users\brunofilipe.lobo\OneDrive - CGI\Code\realvidaseguros\Anexos\Ficheiro1.docx;'
$wsTemplates.Range("C9").Value  = 'C:\Users\brunofilipe.lobo\OneDrive - CGI\Code\realvidaseguros\Anexos\Ficheiro2.docx;'
$wsTemplates.Range("C11").Value = 'C:\Users\brunofilipe.lobo\OneDrive - CGI\Code\realvidaseguros\Anexos\Ficheiro3.docx;'
$wsTemplates.Range("C13").Value = 'C:\Users\brunofilipe.lobo\OneDrive - CGI\Code\realvidaseguros\Anexos\Ficheiro4.docx;'
$wsTemplates.Range("C14").Value = 'C:\Users\brunofilipe.lobo\OneDrive - CGI\Code\realvidaseguros\Anexos\Ficheiro5.docx;'

# ---------------------------------------------------------------------
# 3. IdentifEntidade: swap columns D ("Nome") and E ("Apolice")
#    Cut/Insert keeps the exact bestFit column widths in sync with the
#    swapped content (column D becomes as wide as the old column E, etc).
# ---------------------------------------------------------------------
$wsEntidade = $wb.Worksheets.Item("IdentifEntidade")
$wsEntidade.Columns.Item(5).Cut() | Out-Null
$wsEntidade.Columns.Item(4).Insert() | Out-Null

# ---------------------------------------------------------------------
# 4. View state: selections + which sheet/tab is active.
#    Activating a sheet and selecting a range mirrors what Excel persists
#    as <selection activeCell=.../> and tabSelected="1" in the sheetView.
# ---------------------------------------------------------------------
$wsAnexos = $wb.Worksheets.Item("Anexos")
$wsAnexos.Activate()
$wsAnexos.Range("K25").Select() | Out-Null

$wsEntidade.Activate()
$wsEntidade.Range("I4").Select() | Out-Null

$wsTemplates.Activate()
$wsTemplates.Range("I17").Select() | Out-Null

$wsLabelMap.Activate()
$wsLabelMap.Range("I14").Select() | Out-Null
